$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns stay as text, matching the
# original inline-string formatting (avoids Excel auto-converting values
# like "1.010" or "26.619.25" into numbers).
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '26.619.25'
$ws.Range("E2").Value = '  -2.51%  '

# Row 3
$ws.Range("D3").Value = '1.815.66'
$ws.Range("E3").Value = '  -1.97%  '

# Row 4
$ws.Range("D4").Value = '1.012'
$ws.Range("E4").Value = '  +1.00%  '

# Row 5
$ws.Range("B5").Value = 'USDC'
$ws.Range("C5").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D5").Value = '1.010'
$ws.Range("E5").Value = '  +0.65%  '

# Row 6
$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").Value = '309.46'
$ws.Range("E6").Value = '  -1.47%  '

# Row 7
$ws.Range("D7").Value = '0.4529'
$ws.Range("E7").Value = '  -1.76%  '

# Row 8
$ws.Range("D8").Value = '0.3654'
$ws.Range("E8").Value = '  -1.50%  '

# Row 9
$ws.Range("D9").Value = '0.07075'
$ws.Range("E9").Value = '  -3.06%  '

# Row 10
$ws.Range("D10").Value = '0.8677'
$ws.Range("E10").Value = '  -1.96%  '

# Row 11
$ws.Range("D11").Value = '0.07791'
$ws.Range("E11").Value = '  -0.36%  '

# Row 12
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.886.12'
$ws.Range("E12").Value = '  +2.06%  '

# Row 13
$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D13").Value = '19.12'
$ws.Range("E13").Value = '  -3.92%  '

# Row 14
$ws.Range("D14").Value = '5.281'
$ws.Range("E14").Value = '  -1.79%  '

# Row 15
$ws.Range("D15").Value = '6.293'
$ws.Range("E15").Value = '  -3.89%  '

# Row 16
$ws.Range("D16").Value = '86.24'
$ws.Range("E16").Value = '  -5.87%  '

# Row 17
$ws.Range("D17").Value = '1.015'
$ws.Range("E17").Value = '  +1.08%  '

# Row 18
$ws.Range("D18").Value = '0.000008620'
$ws.Range("E18").Value = '  -3.84%  '

# Row 19
$ws.Range("D19").Value = '1.009'
$ws.Range("E19").Value = '  +0.65%  '

# Row 20
$ws.Range("D20").Value = '26.666.39'
$ws.Range("E20").Value = '  -2.41%  '

# Row 21
$ws.Range("D21").Value = '14.32'
$ws.Range("E21").Value = '  -2.98%  '

# Row 22
$ws.Range("D22").Value = '4.942'
$ws.Range("E22").Value = '  -3.50%  '

# Row 23
$ws.Range("D23").Value = '2.084.32'
$ws.Range("E23").Value = '  +0.22%  '

# Row 24
$ws.Range("D24").Value = '10.35'
$ws.Range("E24").Value = '  -1.71%  '

# Row 25
$ws.Range("D25").Value = '1.989'
$ws.Range("E25").Value = '  +3.12%  '

# Row 26
$ws.Range("D26").Value = '150.83'
$ws.Range("E26").Value = '  -0.65%  '

# Row 27
$ws.Range("D27").Value = '17.98'
$ws.Range("E27").Value = '  -2.30%  '

# Row 28
$ws.Range("D28").Value = '1.981'
$ws.Range("E28").Value = '  -3.33%  '

# Row 29
$ws.Range("D29").Value = '112.91'
$ws.Range("E29").Value = '  -2.65%  '

# Row 30
$ws.Range("D30").Value = '4.860'
$ws.Range("E30").Value = '  -4.23%  '

# Row 31
$ws.Range("D31").Value = '0.08693'
$ws.Range("E31").Value = '  -1.72%  '

# Row 32
$ws.Range("D32").Value = '3.042'
$ws.Range("E32").Value = '  -1.35%  '

# Row 33
$ws.Range("D33").Value = '0.7275'
$ws.Range("E33").Value = '  -5.47%  '

# Row 34
$ws.Range("D34").Value = '4.414'
$ws.Range("E34").Value = '  -1.80%  '

# Row 35
$ws.Range("D35").Value = '1.103'
$ws.Range("E35").Value = '  -5.86%  '

# Row 36
$ws.Range("B36").Value = 'TrustWalletToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D36").Value = '1.076'
$ws.Range("E36").Value = '  -0.17%  '

# Row 37
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").Value = '2.443'
$ws.Range("E37").Value = '  -7.86%  '

# Row 38
$ws.Range("D38").Value = '0.01910'
$ws.Range("E38").Value = '  -2.33%  '

# Row 39
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '2.890'
$ws.Range("E39").Value = '  -2.34%  '

# Row 40
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").Value = '0.05071'
$ws.Range("E40").Value = '  -3.07%  '

# Row 41
$ws.Range("D41").Value = '6.872'
$ws.Range("E41").Value = '  -2.12%  '

# Row 42
$ws.Range("D42").Value = '0.4872'
$ws.Range("E42").Value = '  -5.08%  '

# Row 43
$ws.Range("D43").Value = '0.1567'
$ws.Range("E43").Value = '  -4.17%  '

# Row 44
$ws.Range("D44").Value = '8.105'
$ws.Range("E44").Value = '  -3.36%  '

# Row 45
$ws.Range("D45").Value = '1.010'
$ws.Range("E45").Value = '  +0.70%  '

# Row 46
$ws.Range("D46").Value = '0.4577'
$ws.Range("E46").Value = '  -4.73%  '

# Row 47
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '10.01'
$ws.Range("E47").Value = '  -2.95%  '

# Row 48
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").Value = '102.29'
$ws.Range("E48").Value = '  -0.43%  '

# Row 49
$ws.Range("D49").Value = '1.575'
$ws.Range("E49").Value = '  -4.54%  '

# Row 50
$ws.Range("D50").Value = '0.06027'
$ws.Range("E50").Value = '  -3.03%  '

# Row 51
$ws.Range("D51").Value = '63.99'
$ws.Range("E51").Value = '  -2.30%  '
